# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-25 18:18:28
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists
# the people/systems that recorded a session, separated by ", ". This edit
# rotates that list one position to the left (the first recorder moves to
# the end of the list) for the specific rows touched by the upstream sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,8,10,11,12,13,14,15,17,18,19,20,21,22,29,30,31,32,33,35,37,38,39,40,41,42,44,45,46,47,48,49,56,57,58,59,60,62,64,65,66,67,68,69,71,72,73,74,75,76,83,84,85,86,87,88,89,90,93,95,96,97,99,102,109,110,111,112,113,114,115,116,119,121,122,123,125,128,135,136,137,138,139,140,141,142,145,147,148,149,151,154)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $current = $cell.Text
    $parts = $current -split ", "
    if ($parts.Length -gt 1) {
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $cell.Value = $rotated -join ", "
    }
}
